$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2"=0.943272594242444; "D2"=0.03870760940133522; "E2"=0.4161238364381639; "F2"=0.8971675119293394; "G2"=0.002450953031429082; "K2"=0.4193065313855016; "L2"=0.1014743244794687; "M2"=0.1850376212670142; "O2"=3.165296745607179
    "B3"=0.9244307122180544; "D3"=0.03698799021369581; "E3"=0.4193255635418902; "F3"=0.8955381632624793; "G3"=0.002453521676862919; "K3"=0.3725909496973259; "L3"=0.09557365558983122; "M3"=0.1795559893245127; "O3"=3.174041480071168
    "B4"=0.9132807804026584; "D4"=0.03592142722550307; "E4"=0.4214071962284116; "F4"=0.8950674147691657; "G4"=0.002455184049081051; "K4"=0.3437399356207607; "L4"=0.09199092158220168; "M4"=0.1762727489467188; "O4"=3.181173666974303
    "B5"=0.9088429278881449; "D5"=0.03548412297441672; "E5"=0.4222846344310827; "F5"=0.8950087838763565; "G5"=0.002455882970570784; "K5"=0.3319414112629602; "L5"=0.09054113946716313; "M5"=0.174955652528002; "O5"=3.184523339511458
    "B6"=0.9081124315833904; "D6"=0.03541134825013614; "E6"=0.4224320947188192; "F6"=0.8950070937166927; "G6"=0.002456000325794779; "K6"=0.329979786614274; "L6"=0.09030102325520062; "M6"=0.1747382120823211; "O6"=3.18510631994593
    "B7"=0.913220500838321; "D7"=0.03591554037333822; "E7"=0.4214189115521401; "F7"=0.8950660847057819; "G7"=0.002455193387829734; "K7"=0.3435809838269108; "L7"=0.09197132786643891; "M7"=0.176254901561208; "O7"=3.181217047196697
    "B8"=0.9366892046107012; "D8"=0.0381169201765772; "E8"=0.4172037947951948; "F8"=0.8964957815383414; "G8"=0.002451821052386126; "K8"=0.4032341042013741; "L8"=0.09943145369030049; "M8"=0.1831304913784173; "O8"=3.167946011560929
    "B9"=0.9860194417835828; "D9"=0.04234806609404984; "E9"=0.4098542642052543; "F9"=0.9035031342133806; "G9"=0.002445881142974397; "K9"=0.5188637703234917; "L9"=0.114377938829918; "M9"=0.1972644140758462; "O9"=3.155914314167887
    "B10"=1.024261314769689; "D10"=0.04540364513427875; "E10"=0.4050100383975517; "F10"=0.9112175817117816; "G10"=0.002441923445010364; "K10"=0.6029734410539334; "L10"=0.125550265173203; "M10"=0.2080417054722901; "O10"=3.15561650873596
    "B11"=1.042089125544521; "D11"=0.04678204418647169; "E11"=0.402926223650077; "F11"=0.9152851590207831; "G11"=0.002440210357059396; "K11"=0.6410501477819253; "L11"=0.1306739552340588; "M11"=0.2130292266956104; "O11"=3.157338579624593
    "B12"=1.048901694308057; "D12"=0.04730232162118142; "E12"=0.4021543225959912; "F12"=0.9169057449199158; "G12"=0.00243957414192507; "K12"=0.6554416754273689; "L12"=0.1326200480260837; "M12"=0.2149299895889598; "O12"=3.158257938644198
    "B13"=1.047431755139542; "D13"=0.04719034622339535; "E13"=0.4023198012028768; "F13"=0.9165531520322503; "G13"=0.002439710607321308; "K13"=0.6523434269076063; "L13"=0.1322006626954817; "M13"=0.2145200900043918; "O13"=3.158048050170294
    "B14"=1.042648368535254; "D14"=0.04682488178123378; "E14"=0.402862374510871; "F14"=0.9154168765685853; "G14"=0.002440157765133312; "K14"=0.6422346982341764; "L14"=0.1308339444365032; "M14"=0.2131853617841628; "O14"=3.157408859028351
    "B15"=1.039726409611319; "D15"=0.04660080313814774; "E15"=0.403196954372115; "F15"=0.9147313304470259; "G15"=0.00244043328745428; "K15"=0.6360392358428157; "L15"=0.1299975510737994; "M15"=0.2123693743626021; "O15"=3.157052143565181
    "B16"=1.02310486592188; "D16"=0.04531332801987276; "E16"=0.4051486288447119; "F16"=0.9109629899116669; "G16"=0.002442037150853633; "K16"=0.6004812552051817; "L16"=0.1252162448461007; "M16"=0.20771745790816; "O16"=3.155541357548657
    "B17"=1.013018242296425; "D17"=0.04452051550021707; "E17"=0.4063765855542965; "F17"=0.9087942196426866; "G17"=0.002443043384220584; "K17"=0.5786196778608712; "L17"=0.1222935984826421; "M17"=0.2048853197605354; "O17"=3.155090398318833
    "B18"=1.007257315753321; "D18"=0.04406342039244748; "E18"=0.4070941567745079; "F18"=0.9075993437629393; "G18"=0.002443630362619989; "K18"=0.5660280838558265; "L18"=0.1206164658939457; "M18"=0.2032643443280477; "O18"=3.155005851411062
    "B19"=1.005313758026091; "D19"=0.04390846932619041; "E19"=0.4073390531377417; "F19"=0.9072038034081373; "G19"=0.002443830517028573; "K19"=0.561761819477482; "L19"=0.1200492893248537; "M19"=0.202716886443767; "O19"=3.155007248829349
    "B20"=1.014087778686246; "D20"=0.04460502476761974; "E20"=0.4062447000062788; "F20"=0.9090196508231259; "G20"=0.002442935418779397; "K20"=0.5809486850058931; "L20"=0.1226043165234501; "M20"=0.2051859789807366; "O20"=3.155120308151083
    "B21"=1.044051698049202; "D21"=0.04693227367130248; "E21"=0.4027025412422116; "F21"=0.9157484493043739; "G21"=0.002440026085501772; "K21"=0.6452046229214545; "L21"=0.131235224267229; "M21"=0.2135770763747189; "O21"=3.157589350857876
    "B22"=1.063993474678767; "D22"=0.04844338966372419; "E22"=0.400487739476894; "F22"=0.9206140598635812; "G22"=0.002438197466591283; "K22"=0.6870401928241279; "L22"=0.1369101547630578; "M22"=0.2191316031180648; "O22"=3.160760800447861
    "B23"=1.053317514136495; "D23"=0.04763779029195092; "E23"=0.4016606652222148; "F23"=0.9179743717712228; "G23"=0.002439166792885695; "K23"=0.6647265649364158; "L23"=0.1338782415607938; "M23"=0.21616063637601; "O23"=3.158925556475708
    "B24"=1.013604122859732; "D24"=0.04456682216536478; "E24"=0.4063042893190607; "F24"=0.9089175714750155; "G24"=0.00244298420360527; "K24"=0.5798958127256242; "L24"=0.1224638311244064; "M24"=0.205050028330021; "O24"=3.15510624169238
    "B25"=0.9723219870603259; "D25"=0.04121268910410691; "E25"=0.4117447350559393; "F25"=0.9011570612602924; "G25"=0.002447416393850156; "K25"=0.487729534042785; "L25"=0.1103007914044127; "M25"=0.1933715138762366; "O25"=3.157669985308104
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
